# Applies the crypto-price/volume refresh described in the commit diff.
# Values that are purely numeric-looking decimal strings (e.g. "607.40") are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the sheet's existing inlineStr/text convention for the Price
# column) instead of silently coercing them to a Double via COM Value
# assignment, which would drop formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '66.789.91'
    'E2' = '  +0.25%  '
    'D3' = '3.548.53'
    'E3' = '  -1.47%  '
    'E4' = '  +0.06%  '
    'D5' = '''607.40'
    'E5' = '  -0.29%  '
    'D6' = '''145.70'
    'E6' = '  -2.04%  '
    'D7' = '3.549.90'
    'E7' = '  -1.44%  '
    'D8' = '''1.00'
    'E8' = '  -0.12%  '
    'D9' = '''0.516'
    'E9' = '  +5.43%  '
    'E10' = '  -2.40%  '
    'D11' = '''7.85'
    'E11' = '  -2.21%  '
    'D12' = '''0.414'
    'E12' = '  -0.28%  '
    'D13' = '4.158.76'
    'E13' = '  -1.21%  '
    'D14' = '''0.0000197'
    'E14' = '  -5.95%  '
    'D15' = '''29.17'
    'E15' = '  -2.26%  '
    'D16' = '3.552.64'
    'E16' = '  -1.40%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '66.721.00'
    'E17' = '  +0.01%  '
    'B18' = 'TRON'
    'C18' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D18' = '''0.117'
    'E18' = '  +0.87%  '
    'E19' = '  -4.08%  '
    'D20' = '''6.24'
    'E20' = '  -1.75%  '
    'D21' = '''14.77'
    'E21' = '  -2.23%  '
    'D22' = '''426.65'
    'E22' = '  -0.19%  '
    'D23' = '''0.600'
    'E23' = '  -2.96%  '
    'D24' = '''77.67'
    'E24' = '  -1.58%  '
    'D25' = '3.697.54'
    'E25' = '  -1.18%  '
    'D26' = '''0.999'
    'E26' = '  -0.10%  '
    'E27' = '  -5.78%  '
    'D28' = '''8.05'
    'E28' = '  -2.99%  '
    'E29' = '  -1.26%  '
    'D30' = '''9.09'
    'E30' = '  -2.90%  '
    'D31' = '''1.00'
    'E31' = '  +0.05%  '
    'B32' = 'RenzoRestakedETH'
    'C32' = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
    'D32' = '3.561.40'
    'E32' = '  -1.00%  '
    'B33' = 'Kaspa'
    'C33' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D33' = '''0.157'
    'E33' = '  -1.75%  '
    'D34' = '''24.50'
    'E34' = '  -3.83%  '
    'E35' = '  -0.03%  '
    'D36' = '''1.36'
    'E36' = '  -7.39%  '
    'D37' = '''7.72'
    'E37' = '  -1.80%  '
    'E38' = '  -2.88%  '
    'D39' = '''177.21'
    'E39' = '  -0.19%  '
    'D40' = '''5.32'
    'E40' = '  -5.87%  '
    'D41' = '''0.0832'
    'E41' = '  -2.88%  '
    'D42' = '''5.04'
    'E42' = '  -3.76%  '
    'D43' = '''0.864'
    'E43' = '  -3.90%  '
    'D44' = '''45.45'
    'E44' = '  -1.80%  '
    'D45' = '''1.79'
    'E45' = '  -6.25%  '
    'E46' = '  +0.20%  '
    'D47' = '''2.42'
    'E47' = '  -6.14%  '
    'D48' = '''23.62'
    'E48' = '  -2.93%  '
    'D49' = '''7.16'
    'E49' = '  -0.45%  '
    'E50' = '  -4.79%  '
    'D51' = '''0.924'
    'E51' = '  -3.16%  '
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}

Write-Host "Applied $($changes.Count) cell updates"
